$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 for "Docentes responsaveis:" value
# (label already present at A12; this adds the previously-missing B/C value row)
$ws.Rows(13).Insert()

# Copy the B:C formatting from the row below (already correctly styled) onto
# the new row 13, then clear A13 (no label cell there) and set the values.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Cells.Item(13, 1).Clear()
$ws.Cells.Item(13, 2).Value = '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Cells.Item(13, 3).Value = '5840942 - Marco Aurélio Kondracki de Alcântara'

# Fix the mis-mapped content further down the form (each of these rows had
# its B/C value pointing at the wrong shared string before).
$ws.Cells.Item(14, 2).Value = 'Introdução e conceitos; identificação do problema: tipos de áreas; legislação e normas; geoindicadores de degradação; técnicas de recuperação de áreas degradadas; implementação de planos de recuperação; monitoramento.'
$ws.Cells.Item(14, 3).Value = 'Introdução e conceitos; identificação do problema: tipos de áreas; legislação e normas; geoindicadores de degradação; técnicas de recuperação de áreas degradadas; implementação de planos de recuperação; monitoramento.'

$ws.Cells.Item(16, 2).Value = 'Degradação e recuperação ambiental; geoindicadores de degradação; legislação e normas aplicadas à recuperação de áreas degradadas; aspectos e níveis de recuperação; tipos de áreas degradadas; técnicas e medidas de recuperação de áreas degradadas; critérios para a seleção de alternativas; implementação de planos de recuperação; monitoramento; exemplos de recuperação de áreas degradadas.'
$ws.Cells.Item(16, 3).Value = 'Degradação e recuperação ambiental; geoindicadores de degradação; legislação e normas aplicadas à recuperação de áreas degradadas; aspectos e níveis de recuperação; tipos de áreas degradadas; técnicas e medidas de recuperação de áreas degradadas; critérios para a seleção de alternativas; implementação de planos de recuperação; monitoramento; exemplos de recuperação de áreas degradadas.'

$ws.Cells.Item(19, 2).Value = 'Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'
$ws.Cells.Item(19, 3).Value = 'Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'

$ws.Cells.Item(20, 2).Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'
$ws.Cells.Item(20, 3).Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'

$ws.Cells.Item(21, 2).Value = 'Provas e/ou exercícios dirigidos.'
$ws.Cells.Item(21, 3).Value = 'Provas e/ou exercícios dirigidos.'

$ws.Cells.Item(22, 2).Value = 'Bibliografia básica:Barrow, C.J. Land Degradation Cambridge University Press, 1991.Berger, A.R. The geoindicator concept and its application: An introduction. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changers in Earth Systems. pp.: 1-14 Balkema, Rotterdam, 1996.Blaikie, P. & Brookfield, H. Land degradation and society. London Methuen, 1987Brunsden, D. and Moore, R. Engineering geomorphology on the coast: lessons from West Dorset. Geomorphology 31: 391-409, 1999.CALIJURI, M.C.; CUNHA, D.G.F. Engenharia Ambiental. Conceitos, Tecnologia e Gestão. Rio de Janeiro, Elsevier, 2013.Dahlberg, A.C. Interpretations of environmetal change and diversity: A critical approach to indications of degradation - The case of Kalakamate, Northeast Botswana. Land degradation & Development, 11: 549-562, 2000.DIAS, L.E; GRIFFTITH,J.J. Conceituação e Caracterização de Áreas Degradadas. In: DIAS, L.E; MELLO, J.W.V (orgs). Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Sociedade Brasileira de Recuperação de Áreas Degradadas, 1998.Duque, M.J.F., Pedroza, J., Ciez, A., Sanz, M.A. & Carrasco, R.M.. A geomorphical design for the rehabilitation of an abandoned sand quarry in central Spain. Landscape and urban planning, 42: 1-14, 1998.GUERRA, A. J. T.; ARAUJO, G., ALMEIDA, J. R. Gestão Ambiental De Áreas Degradadas. Rio de Janeiro : Bertrand Brasil, 2007.Marchetti, M. & Panizza, M. Geomorphology and Environmental Impact Assesssment: A case study in Moema (Dolomites - Italy). In: Marchetti, M & Pinas, V. (EDS). Geomorphology and Environemental Impact Assessements pp: 71-82, Balkema, 2001.MARTINS, S.V. Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Viçosa, 2013.Neimanis, U. & kerr, A. Developing national environmental indicators. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing rapid environmental geoindicators: changes in earth systems. 1996.SANCHEZ, L.E. Desengenharia: o passive ambiental na desativaçao de empreendimentos industriais. São Paulo, EDUSP, 2001.SANCHEZ, L.E. Avaliação de Impacto Ambiental. São Paulo, Oficina de Textos, 2006.Bibliografia complementar:Berger, A.R. Assessing Rapid Environmetal Change Using Geoindicators. Environmetal Geology, 32, n. 1, 36-44, 1997.Fao. A provisional methodology for soil degradation assessment. FAO. Rome, 1979,Lindskog, P. and Tengberg, A. Land degradation, Natural resources and local knowledge in the Sahel zone of Burkina Faso. Geojournal, 33, 365-375, 1994.Morton, R. A. Geoindicators of coastal wet land and shorelines. In: berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changes inEarth Systems. pp: 207-232, 1996.Murthy, R.C. Rao, Y. R. and Inamdar, A.B. Integrated coastal management of Mumbai Metropolitan Region. Ocean & Coastal Management 44: 355-369, 2001.'
$ws.Cells.Item(22, 3).Value = 'Bibliografia básica:Barrow, C.J. Land Degradation Cambridge University Press, 1991.Berger, A.R. The geoindicator concept and its application: An introduction. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changers in Earth Systems. pp.: 1-14 Balkema, Rotterdam, 1996.Blaikie, P. & Brookfield, H. Land degradation and society. London Methuen, 1987Brunsden, D. and Moore, R. Engineering geomorphology on the coast: lessons from West Dorset. Geomorphology 31: 391-409, 1999.CALIJURI, M.C.; CUNHA, D.G.F. Engenharia Ambiental. Conceitos, Tecnologia e Gestão. Rio de Janeiro, Elsevier, 2013.Dahlberg, A.C. Interpretations of environmetal change and diversity: A critical approach to indications of degradation - The case of Kalakamate, Northeast Botswana. Land degradation & Development, 11: 549-562, 2000.DIAS, L.E; GRIFFTITH,J.J. Conceituação e Caracterização de Áreas Degradadas. In: DIAS, L.E; MELLO, J.W.V (orgs). Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Sociedade Brasileira de Recuperação de Áreas Degradadas, 1998.Duque, M.J.F., Pedroza, J., Ciez, A., Sanz, M.A. & Carrasco, R.M.. A geomorphical design for the rehabilitation of an abandoned sand quarry in central Spain. Landscape and urban planning, 42: 1-14, 1998.GUERRA, A. J. T.; ARAUJO, G., ALMEIDA, J. R. Gestão Ambiental De Áreas Degradadas. Rio de Janeiro : Bertrand Brasil, 2007.Marchetti, M. & Panizza, M. Geomorphology and Environmental Impact Assesssment: A case study in Moema (Dolomites - Italy). In: Marchetti, M & Pinas, V. (EDS). Geomorphology and Environemental Impact Assessements pp: 71-82, Balkema, 2001.MARTINS, S.V. Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Viçosa, 2013.Neimanis, U. & kerr, A. Developing national environmental indicators. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing rapid environmental geoindicators: changes in earth systems. 1996.SANCHEZ, L.E. Desengenharia: o passive ambiental na desativaçao de empreendimentos industriais. São Paulo, EDUSP, 2001.SANCHEZ, L.E. Avaliação de Impacto Ambiental. São Paulo, Oficina de Textos, 2006.Bibliografia complementar:Berger, A.R. Assessing Rapid Environmetal Change Using Geoindicators. Environmetal Geology, 32, n. 1, 36-44, 1997.Fao. A provisional methodology for soil degradation assessment. FAO. Rome, 1979,Lindskog, P. and Tengberg, A. Land degradation, Natural resources and local knowledge in the Sahel zone of Burkina Faso. Geojournal, 33, 365-375, 1994.Morton, R. A. Geoindicators of coastal wet land and shorelines. In: berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changes inEarth Systems. pp: 207-232, 1996.Murthy, R.C. Rao, Y. R. and Inamdar, A.B. Integrated coastal management of Mumbai Metropolitan Region. Ocean & Coastal Management 44: 355-369, 2001.'
